# Apply "#5: fund, bonds, otherbonds, antique done" edit to the
# "其他有價證券" (Other securities) sheet: add metadata columns
# (property_category .. index) to the sheet, matching the pattern
# already used on the other per-item sheets (stock, deposit, ...),
# and turn row 1 into a proper header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)   # 其他有價證券

# --- Row 1: header labels (was a stray "example" data row) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "quantity"
$ws.Range("E1").Value = "face_value"
$ws.Range("F1").Value = "currency"
$ws.Range("G1").Value = "total"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 ---
$ws.Range("D2").Value = 1192
$ws.Range("H2").Value = "otherbonds"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-03-06"
$ws.Range("K2").Value = "葉宜津"
$ws.Range("L2").Value = 855
$ws.Range("M2").Value = "tmp94f81"
$ws.Range("N2").Value = 110

# --- Row 3 ---
$ws.Range("H3").Value = "otherbonds"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2012-03-06"
$ws.Range("K3").Value = "葉宜津"
$ws.Range("L3").Value = 855
$ws.Range("M3").Value = "tmp94f81"
$ws.Range("N3").Value = 111
